# Update processed/analyzed data for rows 2, 13, and 18 (columns I:L)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = -0.6803438013925298
$ws.Range("J2").Value = 0.1810347255599485
$ws.Range("K2").Value = -0.1615533264691585
$ws.Range("L2").Value = 2.754507615626635

$ws.Range("I13").Value = -0.4897537921999018
$ws.Range("J13").Value = 0.09389166313028015
$ws.Range("K13").Value = 0.4678594645122884
$ws.Range("L13").Value = 2.483941302951547

$ws.Range("I18").Value = -0.8889696259588566
$ws.Range("J18").Value = 0.2256804261455844
$ws.Range("K18").Value = 0.2480950099360809
$ws.Range("L18").Value = 2.281291930843965
